$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before column DW (127) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns.Item(127).Insert()

$ws1.Range("DW1").Value = "22-nov"

for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 127).Value = "-"
}

# --- Sheet "Gaz": append new row 156 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A156").NumberFormat = "@"
$ws2.Range("A156").Value = "2025-11-20"
$ws2.Range("A156").ClearFormats()
$ws2.Range("B156").Value = 30.04

# --- Sheet "CO2": append new row 156 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A156").NumberFormat = "@"
$ws3.Range("A156").Value = "2025-11-20"
$ws3.Range("A156").ClearFormats()
$ws3.Range("B156").Value = 80.92
